$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-06 01:59:58"

# Insert a new row above row 15, shifting existing rows 15-18 down to 16-19.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new job posting.
$ws.Cells.Item(15, 1).Value = $newTimestamp
$ws.Cells.Item(15, 2).Value = "《長期レギュラー》公的機関Web運用の要となる、ディレクター募集"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5465685"
$ws.Cells.Item(15, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://www.lancers.jp/work/detail/5465685")
$ws.Cells.Item(15, 7).Value = 18

# Refresh the "取得日時" timestamp on every other data row (2-14 were pushed
# down? No - rows 2-14 stay put, rows that used to be 15-18 are now 16-19).
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
